# Add a new "Surgery" reference sheet (code / surgery_name / rachs_score)
# right after "PassportPriority", matching the HCA category-codes workbook
# used to import surgical encounters (surgeon/facility/type placeholders
# live on other existing sheets already).

$wb = $excel.ActiveWorkbook

$passportPriority = $wb.Worksheets.Item("PassportPriority")
$surgery = $wb.Worksheets.Add($null, $passportPriority)
$surgery.Name = "Surgery"

# Header row
$surgery.Range("A1").Value = "id"
$surgery.Range("B1").Value = "code"
$surgery.Range("C1").Value = "surgery_name"
$surgery.Range("D1").Value = "rachs_score"

# Fill surgery_name column first, then code column, then id/rachs_score --
# matches the shared-string insertion order of the source workbook.
$surgery.Range("C2").Value = "TOF surgery"
$surgery.Range("C3").Value = "VSD surgery"
$surgery.Range("C4").Value = "ASD surgery"

$surgery.Range("B2").Value = "TOF"
$surgery.Range("B3").Value = "VSD"
$surgery.Range("B4").Value = "ASD"

$surgery.Range("A2").Value = 1
$surgery.Range("D2").Value = 3

$surgery.Range("A3").Value = 2
$surgery.Range("D3").Value = 2

$surgery.Range("A4").Value = 3
$surgery.Range("D4").Value = 2

$surgery.Columns("C").ColumnWidth = 15.5

# The PassportPriority tab is no longer the active/selected one; the user
# clicked an empty cell on it before moving to the new sheet.
[void]$passportPriority.Range("E26").Select()

# Land on the new Surgery sheet, which becomes the active tab.
[void]$surgery.Activate()
[void]$surgery.Range("D5").Select()

Write-Host "Added Surgery sheet with 3 placeholder rows"
